$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price/volume snapshot (and the Stellar/WEMIXTOKEN row swap at 35-36)
$updates = @(
    @{ Cell = 'D2'; Value = '24.789.11' },
    @{ Cell = 'E2'; Value = '  +1.67%  ' },
    @{ Cell = 'D3'; Value = '1.701.19' },
    @{ Cell = 'E3'; Value = '  +1.52%  ' },
    @{ Cell = 'E4'; Value = '  +0.27%  ' },
    @{ Cell = 'D5'; Value = '311.22' },
    @{ Cell = 'E5'; Value = '  +1.97%  ' },
    @{ Cell = 'D6'; Value = '1.001' },
    @{ Cell = 'E6'; Value = '  +0.27%  ' },
    @{ Cell = 'D7'; Value = '0.3725' },
    @{ Cell = 'E7'; Value = '  +1.07%  ' },
    @{ Cell = 'D8'; Value = '49.16' },
    @{ Cell = 'E8'; Value = '  +3.33%  ' },
    @{ Cell = 'D9'; Value = '0.3417' },
    @{ Cell = 'E9'; Value = '  -0.17%  ' },
    @{ Cell = 'E10'; Value = '  +4.44%  ' },
    @{ Cell = 'D11'; Value = '0.07474' },
    @{ Cell = 'E11'; Value = '  +3.72%  ' },
    @{ Cell = 'E12'; Value = '  +0.11%  ' },
    @{ Cell = 'D13'; Value = '20.94' },
    @{ Cell = 'E13'; Value = '  +4.24%  ' },
    @{ Cell = 'D14'; Value = '6.292' },
    @{ Cell = 'E14'; Value = '  +2.69%  ' },
    @{ Cell = 'D15'; Value = '6.999' },
    @{ Cell = 'E15'; Value = '  +4.05%  ' },
    @{ Cell = 'D16'; Value = '1.697.38' },
    @{ Cell = 'E16'; Value = '  +1.32%  ' },
    @{ Cell = 'D17'; Value = '0.00001124' },
    @{ Cell = 'E17'; Value = '  +2.07%  ' },
    @{ Cell = 'D18'; Value = '0.06705' },
    @{ Cell = 'E18'; Value = '  +0.89%  ' },
    @{ Cell = 'E19'; Value = '  +0.22%  ' },
    @{ Cell = 'D20'; Value = '83.37' },
    @{ Cell = 'E20'; Value = '  +3.81%  ' },
    @{ Cell = 'D21'; Value = '17.17' },
    @{ Cell = 'E21'; Value = '  +4.32%  ' },
    @{ Cell = 'D22'; Value = '6.320' },
    @{ Cell = 'E22'; Value = '  +3.68%  ' },
    @{ Cell = 'D23'; Value = '12.91' },
    @{ Cell = 'E23'; Value = '  +6.07%  ' },
    @{ Cell = 'D24'; Value = '24.783.59' },
    @{ Cell = 'E24'; Value = '  +1.92%  ' },
    @{ Cell = 'D25'; Value = '2.454' },
    @{ Cell = 'E25'; Value = '  +0.48%  ' },
    @{ Cell = 'D26'; Value = '2.761' },
    @{ Cell = 'E26'; Value = '  +4.12%  ' },
    @{ Cell = 'D27'; Value = '20.21' },
    @{ Cell = 'E27'; Value = '  +4.32%  ' },
    @{ Cell = 'D28'; Value = '149.20' },
    @{ Cell = 'E28'; Value = '  -2.16%  ' },
    @{ Cell = 'D29'; Value = '131.42' },
    @{ Cell = 'E29'; Value = '  +3.16%  ' },
    @{ Cell = 'D30'; Value = '1.888.24' },
    @{ Cell = 'E30'; Value = '  +1.46%  ' },
    @{ Cell = 'D31'; Value = '1.235' },
    @{ Cell = 'E31'; Value = '  +27.42%  ' },
    @{ Cell = 'D32'; Value = '6.721' },
    @{ Cell = 'E32'; Value = '  +7.08%  ' },
    @{ Cell = 'D33'; Value = '4.226' },
    @{ Cell = 'E33'; Value = '  +4.40%  ' },
    @{ Cell = 'E34'; Value = '  +10.33%  ' },
    @{ Cell = 'B35'; Value = 'Stellar' },
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Cell = 'D35'; Value = '0.08722' },
    @{ Cell = 'E35'; Value = '  +3.14%  ' },
    @{ Cell = 'B36'; Value = 'WEMIXTOKEN' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Cell = 'D36'; Value = '1.762' },
    @{ Cell = 'E36'; Value = '  +5.36%  ' },
    @{ Cell = 'D37'; Value = '5.552' },
    @{ Cell = 'E37'; Value = '  +4.82%  ' },
    @{ Cell = 'D38'; Value = '0.06623' },
    @{ Cell = 'E38'; Value = '  +3.77%  ' },
    @{ Cell = 'D39'; Value = '9.057' },
    @{ Cell = 'E39'; Value = '  +4.56%  ' },
    @{ Cell = 'D40'; Value = '0.02402' },
    @{ Cell = 'E40'; Value = '  +3.89%  ' },
    @{ Cell = 'D41'; Value = '0.2212' },
    @{ Cell = 'E41'; Value = '  +6.21%  ' },
    @{ Cell = 'E42'; Value = '  +3.26%  ' },
    @{ Cell = 'D43'; Value = '0.6399' },
    @{ Cell = 'E43'; Value = '  +5.37%  ' },
    @{ Cell = 'D44'; Value = '1.001' },
    @{ Cell = 'E44'; Value = '  +0.23%  ' },
    @{ Cell = 'D45'; Value = '13.78' },
    @{ Cell = 'E45'; Value = '  +6.46%  ' },
    @{ Cell = 'D46'; Value = '0.6103' },
    @{ Cell = 'E46'; Value = '  +4.09%  ' },
    @{ Cell = 'D47'; Value = '3.817' },
    @{ Cell = 'E47'; Value = '  +1.82%  ' },
    @{ Cell = 'D48'; Value = '2.104' },
    @{ Cell = 'E48'; Value = '  +4.54%  ' },
    @{ Cell = 'D49'; Value = '129.10' },
    @{ Cell = 'E49'; Value = '  +2.84%  ' },
    @{ Cell = 'D50'; Value = '0.07269' },
    @{ Cell = 'E50'; Value = '  +1.77%  ' },
    @{ Cell = 'D51'; Value = '79.43' },
    @{ Cell = 'E51'; Value = '  +5.04%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell -like "D*") {
        # Column D holds text-formatted price strings (e.g. "1.702.20"); force
        # text storage so Excel does not reinterpret them as numbers, then drop
        # the temporary number-format override so the cell keeps its original
        # (unstyled) appearance.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}